$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 Col 1: "68 x 50" -> "51 x 59"
$cell = $t.Cell(1, 1)
$rng = $cell.Range
$rng.End = $rng.End - 2
$rng.Text = "51 x 59" + [char]11 + "  5    9" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "1|    |"

# Row 1 Col 2: "78 x 39" -> "66 x 25"
$cell = $t.Cell(1, 2)
$rng = $cell.Range
$rng.End = $rng.End - 2
$rng.Text = "66 x 25" + [char]11 + "  2    5" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "6|    |"

# Row 1 Col 3: "57 x 48" -> "26 x 98"
$cell = $t.Cell(1, 3)
$rng = $cell.Range
$rng.End = $rng.End - 2
$rng.Text = "26 x 98" + [char]11 + "  9    8" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "6|    |"

# Row 2 Col 1: "10 x 60" -> "40 x 43"
$cell = $t.Cell(2, 1)
$rng = $cell.Range
$rng.End = $rng.End - 2
$rng.Text = "40 x 43" + [char]11 + "  4    3" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "0|    |"

# Row 2 Col 2: "71 x 52" -> "89 x 31"
$cell = $t.Cell(2, 2)
$rng = $cell.Range
$rng.End = $rng.End - 2
$rng.Text = "89 x 31" + [char]11 + "  3    1" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "9|    |"

# Row 2 Col 3: "21 x 97" -> "98 x 60"
$cell = $t.Cell(2, 3)
$rng = $cell.Range
$rng.End = $rng.End - 2
$rng.Text = "98 x 60" + [char]11 + "  6    0" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "8|    |"

# Row 3 Col 1: "51 x 13" -> "35 x 52"
$cell = $t.Cell(3, 1)
$rng = $cell.Range
$rng.End = $rng.End - 2
$rng.Text = "35 x 52" + [char]11 + "  5    2" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "5|    |"

# Row 3 Col 2: "34 x 78" -> "90 x 26"
$cell = $t.Cell(3, 2)
$rng = $cell.Range
$rng.End = $rng.End - 2
$rng.Text = "90 x 26" + [char]11 + "  2    6" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "0|    |"

# Row 3 Col 3: "35 x 14" -> "92 x 83"
$cell = $t.Cell(3, 3)
$rng = $cell.Range
$rng.End = $rng.End - 2
$rng.Text = "92 x 83" + [char]11 + "  8    3" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "2|    |"

# Row 4 Col 1: "32 x 48" -> "56 x 57"
$cell = $t.Cell(4, 1)
$rng = $cell.Range
$rng.End = $rng.End - 2
$rng.Text = "56 x 57" + [char]11 + "  5    7" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "6|    |"

# Row 4 Col 2: "49 x 23" -> "61 x 62"
$cell = $t.Cell(4, 2)
$rng = $cell.Range
$rng.End = $rng.End - 2
$rng.Text = "61 x 62" + [char]11 + "  6    2" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "1|    |"

# Row 4 Col 3: "31 x 51" -> "86 x 86"
$cell = $t.Cell(4, 3)
$rng = $cell.Range
$rng.End = $rng.End - 2
$rng.Text = "86 x 86" + [char]11 + "  8    6" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "6|    |"

# Row 5 Col 1: "46 x 95" -> "57 x 59"
$cell = $t.Cell(5, 1)
$rng = $cell.Range
$rng.End = $rng.End - 2
$rng.Text = "57 x 59" + [char]11 + "  5    9" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "7|    |"

# Row 5 Col 2: "81 x 46" -> "81 x 87"
$cell = $t.Cell(5, 2)
$rng = $cell.Range
$rng.End = $rng.End - 2
$rng.Text = "81 x 87" + [char]11 + "  8    7" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "1|    |"

# Row 5 Col 3: "73 x 96" -> "38 x 89"
$cell = $t.Cell(5, 3)
$rng = $cell.Range
$rng.End = $rng.End - 2
$rng.Text = "38 x 89" + [char]11 + "  8    9" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "8|    |"

